$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.436.92"
$ws.Range("E2").Value = "'  -0.28%  "
$ws.Range("D3").Value = "'1.586.00"
$ws.Range("E3").Value = "'  -0.07%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'213.60"
$ws.Range("E6").Value = "'  -0.26%  "
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("D8").Value = "'44.51"
$ws.Range("E8").Value = "'  -0.20%  "
$ws.Range("D9").Value = "'23.96"
$ws.Range("E9").Value = "'  -1.15%  "
$ws.Range("E10").Value = "'  -1.80%  "
$ws.Range("E11").Value = "'  -1.61%  "
$ws.Range("E12").Value = "'  +0.94%  "
$ws.Range("D13").Value = "'1.812.27"
$ws.Range("E13").Value = "'  -0.08%  "
$ws.Range("D14").Value = "'1.581.63"
$ws.Range("E14").Value = "'  -0.30%  "
$ws.Range("D15").Value = "'3.71"
$ws.Range("E15").Value = "'  -0.74%  "
$ws.Range("E16").Value = "'  -1.53%  "
$ws.Range("D17").Value = "'28.473.97"
$ws.Range("E17").Value = "'  -0.20%  "
$ws.Range("D18").Value = "'62.17"
$ws.Range("E18").Value = "'  -1.43%  "
$ws.Range("D19").Value = "'230.36"
$ws.Range("E19").Value = "'  -0.22%  "
$ws.Range("D20").Value = "'7.46"
$ws.Range("E20").Value = "'  -0.51%  "
$ws.Range("D21").Value = "'0.0₃0690"
$ws.Range("E21").Value = "'  -2.34%  "
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("D23").Value = "'3.92"
$ws.Range("E23").Value = "'  -3.12%  "
$ws.Range("D24").Value = "'9.16"
$ws.Range("E24").Value = "'  -1.77%  "
$ws.Range("E25").Value = "'  +3.40%  "
$ws.Range("D26").Value = "'151.80"
$ws.Range("E26").Value = "'  -0.18%  "
$ws.Range("D27").Value = "'15.02"
$ws.Range("E27").Value = "'  -1.29%  "
$ws.Range("E28").Value = "'  -1.64%  "
$ws.Range("E29").Value = "'  -1.86%  "
$ws.Range("E30").Value = "'  -0.01%  "
$ws.Range("E32").Value = "'  -1.08%  "
$ws.Range("E33").Value = "'  -1.16%  "
$ws.Range("E34").Value = "'  -2.27%  "
$ws.Range("D35").Value = "'1.395.85"
$ws.Range("E35").Value = "'  +0.57%  "
$ws.Range("D37").Value = "'1.51"
$ws.Range("E37").Value = "'  -4.60%  "
$ws.Range("D38").Value = "'2.36"
$ws.Range("E38").Value = "'  +0.38%  "
$ws.Range("E39").Value = "'  +1.43%  "
$ws.Range("E40").Value = "'  -0.70%  "
$ws.Range("D41").Value = "'0.523"
$ws.Range("E42").Value = "'  +0.00%  "
$ws.Range("D43").Value = "'0.792"
$ws.Range("E43").Value = "'  -2.46%  "
$ws.Range("B45").Value = "'FraxShare"
$ws.Range("C45").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.45"
$ws.Range("E45").Value = "'  -3.27%  "
$ws.Range("B46").Value = "'Kaspa"
$ws.Range("C46").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").Value = "'0.0458"
$ws.Range("E46").Value = "'  -1.44%  "
$ws.Range("E47").Value = "'  -2.07%  "
$ws.Range("D48").Value = "'63.02"
$ws.Range("E48").Value = "'  +0.17%  "
$ws.Range("D49").Value = "'1.723.21"
$ws.Range("E49").Value = "'  -0.02%  "
$ws.Range("D50").Value = "'86.64"
$ws.Range("E50").Value = "'  -0.37%  "
$ws.Range("D51").Value = "'0.0₆0103"
$ws.Range("E51").Value = "'  -2.45%  "
